$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (e.g. "1.008", "27.819.98") are written
# as literal text, matching the source data feed which stores these as strings.
$textCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "D26",
    "E26",
    "D27",
    "E27",
    "D28",
    "E28",
    "D29",
    "E29",
    "D30",
    "E30",
    "D31",
    "E31",
    "D32",
    "E32",
    "D33",
    "E33",
    "D34",
    "E34",
    "D35",
    "E35",
    "D36",
    "E36",
    "D37",
    "E37",
    "D38",
    "E38",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '27.819.98'
$ws.Range("E2").Value = '  +1.17%  '

# Row 3
$ws.Range("D3").Value = '1.901.07'
$ws.Range("E3").Value = '  +2.22%  '

# Row 4
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  -0.26%  '

# Row 5
$ws.Range("D5").Value = '315.96'
$ws.Range("E5").Value = '  +1.50%  '

# Row 6
$ws.Range("E6").Value = '  -0.27%  '

# Row 7
$ws.Range("D7").Value = '0.4827'
$ws.Range("E7").Value = '  +1.31%  '

# Row 8
$ws.Range("D8").Value = '0.3806'
$ws.Range("E8").Value = '  +0.39%  '

# Row 9
$ws.Range("D9").Value = '0.07350'
$ws.Range("E9").Value = '  +0.34%  '

# Row 10
$ws.Range("D10").Value = '0.9266'
$ws.Range("E10").Value = '  -0.25%  '

# Row 11
$ws.Range("D11").Value = '20.66'
$ws.Range("E11").Value = '  -0.13%  '

# Row 12
$ws.Range("D12").Value = '0.07737'
$ws.Range("E12").Value = '  -0.45%  '

# Row 13
$ws.Range("D13").Value = '1.924.76'
$ws.Range("E13").Value = '  +2.91%  '

# Row 14
$ws.Range("D14").Value = '5.478'
$ws.Range("E14").Value = '  +0.75%  '

# Row 15
$ws.Range("D15").Value = '6.601'
$ws.Range("E15").Value = '  +0.60%  '

# Row 16
$ws.Range("D16").Value = '91.46'
$ws.Range("E16").Value = '  +1.46%  '

# Row 18
$ws.Range("D18").Value = '0.000008842'
$ws.Range("E18").Value = '  +0.34%  '

# Row 19
$ws.Range("D19").Value = '1.006'
$ws.Range("E19").Value = '  -0.27%  '

# Row 20
$ws.Range("D20").Value = '27.887.19'
$ws.Range("E20").Value = '  +1.33%  '

# Row 21
$ws.Range("D21").Value = '14.62'
$ws.Range("E21").Value = '  -0.22%  '

# Row 22
$ws.Range("D22").Value = '5.144'
$ws.Range("E22").Value = '  +1.12%  '

# Row 23
$ws.Range("D23").Value = '2.135.09'
$ws.Range("E23").Value = '  +1.65%  '

# Row 24
$ws.Range("D24").Value = '10.84'
$ws.Range("E24").Value = '  +1.26%  '

# Row 25
$ws.Range("D25").Value = '1.918'
$ws.Range("E25").Value = '  -0.92%  '

# Row 26
$ws.Range("D26").Value = '154.49'
$ws.Range("E26").Value = '  -0.78%  '

# Row 27
$ws.Range("D27").Value = '18.44'
$ws.Range("E27").Value = '  -0.02%  '

# Row 28
$ws.Range("D28").Value = '2.119'
$ws.Range("E28").Value = '  +5.81%  '

# Row 29
$ws.Range("D29").Value = '117.01'
$ws.Range("E29").Value = '  +1.49%  '

# Row 30
$ws.Range("D30").Value = '4.941'
$ws.Range("E30").Value = '  -0.21%  '

# Row 31
$ws.Range("D31").Value = '0.08960'
$ws.Range("E31").Value = '  +1.06%  '

# Row 32
$ws.Range("D32").Value = '3.229'
$ws.Range("E32").Value = '  -3.01%  '

# Row 33
$ws.Range("D33").Value = '1.246'
$ws.Range("E33").Value = '  +3.76%  '

# Row 34
$ws.Range("D34").Value = '0.7629'
$ws.Range("E34").Value = '  +1.27%  '

# Row 35
$ws.Range("D35").Value = '4.642'
$ws.Range("E35").Value = '  +1.30%  '

# Row 36
$ws.Range("D36").Value = '0.02039'
$ws.Range("E36").Value = '  -0.14%  '

# Row 37
$ws.Range("D37").Value = '2.513'
$ws.Range("E37").Value = '  -7.23%  '

# Row 38
$ws.Range("D38").Value = '1.094'
$ws.Range("E38").Value = '  -2.29%  '

# Row 39
$ws.Range("D39").Value = '0.05273'
$ws.Range("E39").Value = '  -1.00%  '

# Row 40
$ws.Range("D40").Value = '2.995'
$ws.Range("E40").Value = '  +0.50%  '

# Row 41
$ws.Range("D41").Value = '0.5456'
$ws.Range("E41").Value = '  -1.87%  '

# Row 42
$ws.Range("D42").Value = '6.941'
$ws.Range("E42").Value = '  -1.27%  '

# Row 43
$ws.Range("D43").Value = '0.1521'
$ws.Range("E43").Value = '  +0.04%  '

# Row 44
$ws.Range("D44").Value = '8.344'
$ws.Range("E44").Value = '  -1.52%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '10.69'
$ws.Range("E45").Value = '  -0.17%  '

# Row 46
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '108.97'
$ws.Range("E46").Value = '  +4.82%  '

# Row 47
$ws.Range("D47").Value = '0.4795'
$ws.Range("E47").Value = '  -1.39%  '

# Row 48
$ws.Range("D48").Value = '1.007'
$ws.Range("E48").Value = '  -0.27%  '

# Row 49
$ws.Range("D49").Value = '1.643'
$ws.Range("E49").Value = '  -1.26%  '

# Row 50
$ws.Range("D50").Value = '67.67'
$ws.Range("E50").Value = '  +0.51%  '

# Row 51
$ws.Range("D51").Value = '0.06083'
$ws.Range("E51").Value = '  -0.29%  '
